$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 68, pushing existing rows 68:150 down to 69:151
$ws.Rows.Item(68).Insert()

# Populate the newly inserted row 68 with the new record
$ws.Range("A68").Value = 10
$ws.Range("B68").Value = "Vega Modelo de Temuco"
$ws.Range("C68").Value = "La Araucanía"
$ws.Range("D68").Value = 44781
$ws.Range("E68").Value = 9
$ws.Range("F68").Value = "Fruta"
$ws.Range("G68").Value = 100104
$ws.Range("H68").Value = "Frutos de pepita"
$ws.Range("I68").Value = 100104001
$ws.Range("J68").Value = "Granada"
$ws.Range("K68").Value = "Wonderfull"
$ws.Range("L68").Value = "Primera"
$ws.Range("M68").Value = 55
$ws.Range("N68").Value = 15000
$ws.Range("O68").Value = 15000
$ws.Range("P68").Value = 15000
$ws.Range("Q68").Value = "$/bandeja 10 kilos granel"
$ws.Range("R68").Value = "Provincia de Limarí"
$ws.Range("S68").Value = 1500
$ws.Range("T68").Value = 10
